$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 577.4400000000001
$ws.Range("I15").Value = 577.4400000000001
$ws.Range("K15").Value = 1732.32
$ws.Range("M15").Value = -1563.32
$ws.Range("H87").Value = 13346.803
$ws.Range("J87").Value = 13346.803
$ws.Range("L87").Value = 13346.803
$ws.Range("N87").Value = -15842.803
$ws.Range("H90").Value = 13346.803
$ws.Range("J90").Value = 13346.803
$ws.Range("L90").Value = 40040.409
$ws.Range("N90").Value = -52520.409
$ws.Range("H98").Value = 2179.1333
$ws.Range("I98").Value = 2386.7778
$ws.Range("J98").Value = 1867.6666
$ws.Range("K98").Value = 2386.7778
$ws.Range("L98").Value = 1867.6666
$ws.Range("M98").Value = -888.7777999999998
$ws.Range("N98").Value = -4863.6666
$ws.Range("H122").Value = 2179.1333
$ws.Range("I122").Value = 2386.7778
$ws.Range("J122").Value = 1867.6666
$ws.Range("K122").Value = 7160.3334
$ws.Range("L122").Value = 5602.9998
$ws.Range("M122").Value = -4710.3334
$ws.Range("N122").Value = -10502.9998
$ws.Range("H132").Value = 25860.334
$ws.Range("I132").Value = 35939.184
$ws.Range("J132").Value = 5702.6313
$ws.Range("K132").Value = 107817.552
$ws.Range("L132").Value = 17107.8939
$ws.Range("M132").Value = -105287.552
$ws.Range("N132").Value = -22167.8939
$ws.Range("H136").Value = 46136
$ws.Range("J136").Value = 46136
$ws.Range("L136").Value = 46136
$ws.Range("N136").Value = -56336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6229.5156
$ws.Range("I32").Value = 5567.785
$ws.Range("K32").Value = 5567.785
$ws.Range("M32").Value = -5280.785
$ws.Range("H74").Value = 1821.5
$ws.Range("I74").Value = 1401.5
$ws.Range("J74").Value = 2381.5
$ws.Range("K74").Value = 1401.5
$ws.Range("L74").Value = 2381.5
$ws.Range("M74").Value = -527.5
$ws.Range("N74").Value = -4129.5
$ws.Range("H77").Value = 1821.5
$ws.Range("I77").Value = 1401.5
$ws.Range("J77").Value = 2381.5
$ws.Range("K77").Value = 7007.5
$ws.Range("L77").Value = 11907.5
$ws.Range("M77").Value = -2639.5
$ws.Range("N77").Value = -20643.5
$ws.Range("H97").Value = 591.21875
$ws.Range("I97").Value = 574.7778
$ws.Range("J97").Value = 680
$ws.Range("K97").Value = 574.7778
$ws.Range("L97").Value = 680
$ws.Range("M97").Value = -78.77779999999996
$ws.Range("N97").Value = -1672
$ws.Range("H132").Value = 1658.4906
$ws.Range("I132").Value = 1253.8684
$ws.Range("J132").Value = 2683.5334
$ws.Range("K132").Value = 3761.6052
$ws.Range("L132").Value = 8050.600199999999
$ws.Range("M132").Value = -1231.6052
$ws.Range("N132").Value = -13110.6002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 35241.668
$ws.Range("J51").Value = 35241.668
$ws.Range("L51").Value = 35241.668
$ws.Range("N51").Value = -36223.668
$ws.Range("H134").Value = 1605.2559
$ws.Range("I134").Value = 1230.963
$ws.Range("J134").Value = 2236.875
$ws.Range("K134").Value = 3692.889
$ws.Range("L134").Value = 6710.625
$ws.Range("M134").Value = -1157.889
$ws.Range("N134").Value = -11780.625
$ws.Range("H140").Value = 50306.363
$ws.Range("J140").Value = 50306.363
$ws.Range("L140").Value = 50306.363
$ws.Range("N140").Value = -60666.363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 28000
$ws.Range("J52").Value = 28000
$ws.Range("L52").Value = 28000
$ws.Range("N52").Value = -28588
$ws.Range("H58").Value = 1936.1632
$ws.Range("I58").Value = 1230.4062
$ws.Range("J58").Value = 3264.647
$ws.Range("K58").Value = 1230.4062
$ws.Range("L58").Value = 3264.647
$ws.Range("M58").Value = -1027.4062
$ws.Range("N58").Value = -3670.647
$ws.Range("H134").Value = 2137.8235
$ws.Range("I134").Value = 2476.9092
$ws.Range("J134").Value = 1516.1666
$ws.Range("K134").Value = 7430.7276
$ws.Range("L134").Value = 4548.4998
$ws.Range("M134").Value = -4895.7276
$ws.Range("N134").Value = -9618.4998
$ws.Range("H136").Value = 1936.1632
$ws.Range("I136").Value = 1230.4062
$ws.Range("J136").Value = 3264.647
$ws.Range("K136").Value = 3691.2186
$ws.Range("L136").Value = 9793.940999999999
$ws.Range("M136").Value = -1141.2186
$ws.Range("N136").Value = -14893.941
$ws.Range("H137").Value = 33698.375
$ws.Range("J137").Value = 33698.375
$ws.Range("L137").Value = 33698.375
$ws.Range("N137").Value = -43898.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 858.04
$ws.Range("J131").Value = 870.15625
$ws.Range("L131").Value = 2610.46875
$ws.Range("N131").Value = -12690.46875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 9998.5
$ws.Range("J86").Value = 9998.5
$ws.Range("L86").Value = 9998.5
$ws.Range("N86").Value = -12370.5
$ws.Range("H89").Value = 9998.5
$ws.Range("J89").Value = 9998.5
$ws.Range("L89").Value = 29995.5
$ws.Range("N89").Value = -41851.5
$ws.Range("H122").Value = 2511.1052
$ws.Range("I122").Value = 2037.7858
$ws.Range("J122").Value = 3836.4
$ws.Range("K122").Value = 6113.357400000001
$ws.Range("L122").Value = 11509.2
$ws.Range("M122").Value = -3663.357400000001
$ws.Range("N122").Value = -16409.2
$ws.Range("H126").Value = 1930.0377
$ws.Range("I126").Value = 1807.7333
$ws.Range("J126").Value = 2089.5652
$ws.Range("K126").Value = 5423.199900000001
$ws.Range("L126").Value = 6268.6956
$ws.Range("M126").Value = -2953.199900000001
$ws.Range("N126").Value = -11208.6956
$ws.Range("H136").Value = 11178.787
$ws.Range("J136").Value = 11178.787
$ws.Range("L136").Value = 33536.361
$ws.Range("N136").Value = -38636.361

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H100").Value = 2535.5806
$ws.Range("I100").Value = 2478.1924
$ws.Range("K100").Value = 2478.1924
$ws.Range("M100").Value = -1937.1924
$ws.Range("H122").Value = 8549168
$ws.Range("I122").Value = 13890076
$ws.Range("J122").Value = 3716
$ws.Range("K122").Value = 41670228
$ws.Range("L122").Value = 11148
$ws.Range("M122").Value = -41667778
$ws.Range("N122").Value = -16048
$ws.Range("H133").Value = 38409.75
$ws.Range("J133").Value = 38409.75
$ws.Range("L133").Value = 38409.75
$ws.Range("N133").Value = -43469.75
$ws.Range("H136").Value = 19611748
$ws.Range("I136").Value = 3989.9092
$ws.Range("J136").Value = 55559304
$ws.Range("K136").Value = 11969.7276
$ws.Range("L136").Value = 166677912
$ws.Range("M136").Value = -9419.7276
$ws.Range("N136").Value = -166683012

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 42485.2
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 42485.2
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 42485.2
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -42947.2
$ws.Range("H119").Value = 15330.077
$ws.Range("J119").Value = 15330.077
$ws.Range("L119").Value = 15330.077
$ws.Range("N119").Value = -25006.077
$ws.Range("H134").Value = 42485.2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 42485.2
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 127455.6
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -132525.6
